$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend a note to the instructions/description cell (A2) to explain the
# file is mechanically generated.
$ws.Range("A2").Value = "Mechanically generated keynote file. REMEMBER TO SAVE after editing, then SAVE FILE AS Text (Tab delimited)(*.txt), then load/reload your keynotes on your project Revit file so Revit can see the changes. All keynote / text editing shall be on the Excel file only."

# The DEMO 1 keynote description (row 6) is long enough to need a taller,
# wrapped row - bump its height so the text is fully visible.
$ws.Rows.Item(6).RowHeight = 30

# Touch every cell covered by a merged range so they're materialised in the
# sheet (matches how Excel writes out merged areas), without altering their
# appearance.
$ws.Range("B2").Borders.LineStyle = -4142
$ws.Range("B10").Borders.LineStyle = -4142
$ws.Range("C10").Borders.LineStyle = -4142
$ws.Range("B14").Borders.LineStyle = -4142
$ws.Range("C14").Borders.LineStyle = -4142
$ws.Range("B18").Borders.LineStyle = -4142
$ws.Range("C18").Borders.LineStyle = -4142
